# Generate Report for Handoff
# Reproduces the localization-status report refresh: status text, handoff/handback
# timestamps, priority, and a new "stale handback" error message surfaced on the
# 9b684054 row for both locales, plus the narrower Status-column / wider
# Error-Detail-column widths that went with the longer error text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("C3").Value = "Ready for handoff"

# --- Latest HO Xliff Generate Date: 23:56:39 -> 23:58:44 ---
$overview.Range("G2").Value = "2016-11-08 23:58:44"
$overview.Range("G3").Value = "2016-11-08 23:58:44"
$dede.Range("H2").Value = "2016-11-08 23:58:44"
$dede.Range("H3").Value = "2016-11-08 23:58:44"

# --- Priority for the 1119d81b file: "ht" -> "mt" ---
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E3").Value = "mt"
$dede.Range("E2").Value = "mt"
$dede.Range("E3").Value = "mt"

# --- Latest Handoff Datetime (zh-cn, 1119d81b row pair): 23:56:20 -> 23:58:30 ---
$zhcn.Range("H2").Value = "2016-11-08 23:58:30"
$zhcn.Range("H3").Value = "2016-11-08 23:58:30"

# --- New Error Detail message for the stale 9b684054 handback, both locales ---
$errMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3b31aa4913947ae6a0eceb9ce7abca696c139a5/e2e/9b684054-780e-4a2c-8e8d-105b67421461.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd67bd9387798966497d77c8d61962854dbfa04b/e2e/9b684054-780e-4a2c-8e8d-105b67421461.md."
$zhcn.Range("P3").Value = $errMsg
$dede.Range("P3").Value = $errMsg

# --- Column width adjustments that accompanied the longer Status / Error Detail text ---
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333

$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
